# Update the cached "today" text shown in the Date placeholders that live on
# the slide master, every slide layout and the notes master (PowerPoint
# refreshes these "Update automatically" fields whenever the deck is
# resaved), and refresh the link-preview text on the "Reading from input
# example" slide to point at the new JSFiddle instead of the old Replit.

$p = $ppt.ActivePresentation

function Update-DatePlaceholders {
    param($container)

    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)

        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "4/1/2024") {
                $tr.Text = "9/4/2024"
            } elseif ($tr.Text -eq "April 1, 2024") {
                $tr.Text = "September 4, 2024"
            }
        }
    }
}

# Slide master.
Update-DatePlaceholders $p.SlideMaster

# Every slide layout hanging off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholders $layouts.Item($j)
}

# Notes master.
Update-DatePlaceholders $p.NotesMaster

# Slide 10 ("Reading from input example") - swap the displayed link text.
$s10 = $p.Slides.Item(10)
for ($k = 1; $k -le $s10.Shapes.Count; $k++) {
    $shp = $s10.Shapes.Item($k)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "https://replit.com/@HylandOutreach/ReadFromDom") {
        $shp.TextFrame.TextRange.Text = "https://jsfiddle.net/d3rme58n/"
    }
}
